$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "24.764.86"
$ws.Range("E2").Value = "  +2.46%  "

# Row 3
$ws.Range("D3").Value = "1.705.44"
$ws.Range("E3").Value = "  +1.53%  "

# Row 4
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  -0.23%  "

# Row 5
$ws.Range("D5").Value = "'309.38"
$ws.Range("E5").Value = "  +0.16%  "

# Row 6
$ws.Range("D6").Value = "'0.9985"
$ws.Range("E6").Value = "  -0.09%  "

# Row 8
$ws.Range("D8").Value = "'49.10"
$ws.Range("E8").Value = "  +3.60%  "

# Row 9
$ws.Range("D9").Value = "'0.3442"
$ws.Range("E9").Value = "  -0.11%  "

# Row 10
$ws.Range("D10").Value = "'1.201"
$ws.Range("E10").Value = "  +0.57%  "

# Row 11
$ws.Range("D11").Value = "'0.07483"
$ws.Range("E11").Value = "  +2.43%  "

# Row 12
$ws.Range("D12").Value = "'0.9986"
$ws.Range("E12").Value = "  -0.27%  "

# Row 13
$ws.Range("D13").Value = "'20.94"
$ws.Range("E13").Value = "  +2.49%  "

# Row 14
$ws.Range("D14").Value = "'6.243"
$ws.Range("E14").Value = "  +2.15%  "

# Row 15
$ws.Range("D15").Value = "'6.979"
$ws.Range("E15").Value = "  +2.99%  "

# Row 16
$ws.Range("D16").Value = "1.707.39"
$ws.Range("E16").Value = "  +1.75%  "

# Row 17
$ws.Range("D17").Value = "'0.00001130"
$ws.Range("E17").Value = "  +2.05%  "

# Row 18
$ws.Range("D18").Value = "'0.06718"
$ws.Range("E18").Value = "  -0.03%  "

# Row 19
$ws.Range("D19").Value = "'0.9983"
$ws.Range("E19").Value = "  -0.13%  "

# Row 20
$ws.Range("D20").Value = "'84.37"
$ws.Range("E20").Value = "  +3.17%  "

# Row 21
$ws.Range("D21").Value = "'17.19"
$ws.Range("E21").Value = "  +4.11%  "

# Row 22
$ws.Range("D22").Value = "'6.334"
$ws.Range("E22").Value = "  +3.63%  "

# Row 23
$ws.Range("D23").Value = "'13.08"
$ws.Range("E23").Value = "  +8.79%  "

# Row 24
$ws.Range("D24").Value = "24.732.20"
$ws.Range("E24").Value = "  +2.45%  "

# Row 25
$ws.Range("D25").Value = "'2.437"
$ws.Range("E25").Value = "  +0.67%  "

# Row 26
$ws.Range("D26").Value = "'2.763"
$ws.Range("E26").Value = "  +3.39%  "

# Row 27
$ws.Range("D27").Value = "'20.29"
$ws.Range("E27").Value = "  +3.37%  "

# Row 28
$ws.Range("D28").Value = "'149.94"
$ws.Range("E28").Value = "  -1.73%  "

# Row 29
$ws.Range("E29").Value = "  +3.38%  "

# Row 30
$ws.Range("D30").Value = "1.894.60"
$ws.Range("E30").Value = "  +1.73%  "

# Row 31
$ws.Range("D31").Value = "'1.181"
$ws.Range("E31").Value = "  +20.80%  "

# Row 32
$ws.Range("D32").Value = "'6.763"
$ws.Range("E32").Value = "  +5.69%  "

# Row 33
$ws.Range("D33").Value = "'4.202"
$ws.Range("E33").Value = "  +3.63%  "

# Row 34
$ws.Range("D34").Value = "'1.794"
$ws.Range("E34").Value = "  +2.21%  "

# Row 35
$ws.Range("D35").Value = "'0.08819"
$ws.Range("E35").Value = "  +4.26%  "

# Row 36
$ws.Range("D36").Value = "'13.66"
$ws.Range("E36").Value = "  +10.68%  "

# Row 37
$ws.Range("D37").Value = "'5.536"
$ws.Range("E37").Value = "  +3.20%  "

# Row 38
$ws.Range("D38").Value = "'0.06599"
$ws.Range("E38").Value = "  +2.21%  "

# Row 39
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.02395"
$ws.Range("E39").Value = "  +1.97%  "

# Row 40
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "'8.982"
$ws.Range("E40").Value = "  +0.31%  "

# Row 41
$ws.Range("D41").Value = "'0.2224"
$ws.Range("E41").Value = "  +4.66%  "

# Row 42
$ws.Range("E42").Value = "  +1.08%  "

# Row 43
$ws.Range("D43").Value = "'0.6452"
$ws.Range("E43").Value = "  +4.28%  "

# Row 44
$ws.Range("D44").Value = "'0.9979"
$ws.Range("E44").Value = "  -0.08%  "

# Row 45
$ws.Range("D45").Value = "'13.93"
$ws.Range("E45").Value = "  +5.42%  "

# Row 46
$ws.Range("D46").Value = "'0.6121"
$ws.Range("E46").Value = "  +2.60%  "

# Row 47
$ws.Range("D47").Value = "'3.816"
$ws.Range("E47").Value = "  +0.31%  "

# Row 48
$ws.Range("D48").Value = "'2.120"
$ws.Range("E48").Value = "  +3.95%  "

# Row 49
$ws.Range("D49").Value = "'129.61"
$ws.Range("E49").Value = "  +2.44%  "

# Row 50
$ws.Range("D50").Value = "'0.07289"
$ws.Range("E50").Value = "  +1.69%  "

# Row 51
$ws.Range("D51").Value = "'79.33"
$ws.Range("E51").Value = "  +4.19%  "
